$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Report")

# Update existing rows: change Status from "In Progress" to "Complete" for rows 9 and 10
$ws.Range("C9").Value = "Complete"
$ws.Range("C10").Value = "Complete"

# Add new row 13
$ws.Range("A13").Value = $ws.Range("A12").Value2
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)

$ws.Range("B13").Value = "Added basic console ui"
$ws.Range("B12").Copy()
$ws.Range("B13").PasteSpecial(-4122)

$ws.Range("C13").Value = "In Progress"
$ws.Range("C12").Copy()
$ws.Range("C13").PasteSpecial(-4122)

$ws.Range("D13").Value = "Caleb Viverito"
$ws.Range("D12").Copy()
$ws.Range("D13").PasteSpecial(-4122)

# Update selection to match final state
$ws.Range("C10").Select()
